$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: "JB" -> "zer"
$ws.Range("E2").Value = "zer"

# Row 18: "lait concentré" (ref 258, 34000/58000) -> "lait en poudre" (ref 40, 20000/44000)
# A18 holds a number-looking label ("40") so it must stay text like the rest
# of the reference column. Write it with a leading apostrophe (forces text)
# then restore the cell's original formatting with a formats-only paste so
# the quote-prefix flag doesn't leave a different style behind.
$ws.Range("A18").Value = "'40"
$ws.Range("B18").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

$ws.Range("B18").Value = "lait en poudre"
$ws.Range("D18").Value = 20000
$ws.Range("G18").Value = 44000

# Rows 19 and 20 (the old "lait en poudre" / "yaourt" lines) are removed and
# become blank rows, matching the style of the blank rows beneath them (s=25).
$ws.Range("A19:G20").ClearContents()
$ws.Range("A21:G21").Copy() | Out-Null
$ws.Range("A19:G20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
